# FedExShipments.xlsx - 28th March 2022 update
# Re-run results for three FedEx shipment rows: tracking numbers and
# actual rates were refreshed, and the PASS/FAIL result recalculated.
# Row 23 and Row 25 now PASS (actual rate matches expected rate); Row 24
# stays FAIL (actual rate still does not match expected rate).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as plain text (matching ShipmentTracking/ActualRate/Result
# columns, which already store text, not numbers/currency) instead of
# letting Excel auto-convert the digit strings / "$"-values to numbers.
# Only touch the specific cells that actually change value below.
$ws.Range("P23:Q24").NumberFormat = "@"
$ws.Range("P25:Q25").NumberFormat = "@"
$ws.Range("R23").NumberFormat = "@"
$ws.Range("R25").NumberFormat = "@"

# Row 23: new tracking number, actual rate now equals expected ($439.28) -> PASS
$ws.Range("P23").Value = "320018126760"
$ws.Range("Q23").Value = "$439.28"
$ws.Range("R23").Value = "PASS"

# Row 24: new tracking number, actual rate still differs from expected -> stays FAIL
$ws.Range("P24").Value = "320018126770"
$ws.Range("Q24").Value = "$278.12"

# Row 25: new tracking number, actual rate now equals expected ($52.88) -> PASS
$ws.Range("P25").Value = "320018126781"
$ws.Range("Q25").Value = "$52.88"
$ws.Range("R25").Value = "PASS"
